$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting old D:K to E:L
$ws.Columns("D").Insert()

# Copy number formatting from column E (shifted original D) into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D (and a few corrected cells) with FY2018 data
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 4642100
$ws.Range("D9").Value2 = 2418500
$ws.Range("D10").Value2 = 2223600
$ws.Range("D12").Value2 = 18500
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = "NA"
$ws.Range("D15").Value2 = "NA"
$ws.Range("D17").Value2 = 4204300
$ws.Range("D18").Value2 = 437800
$ws.Range("D20").Value2 = 0
$ws.Range("D21").Value2 = 547400
$ws.Range("D22").Value2 = 5800
$ws.Range("D23").Value2 = 431900
$ws.Range("D24").Value2 = 71600
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 360300
$ws.Range("D27").Value2 = 290100
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 11000
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 0
$ws.Range("D33").Value2 = 301000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 301000
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 872200
$ws.Range("D42").Value2 = 100000
$ws.Range("E42").Value2 = "NA"
$ws.Range("F42").Value2 = "NA"
$ws.Range("G42").Value2 = "NA"
$ws.Range("H42").Value2 = "NA"
$ws.Range("I42").Value2 = "NA"
$ws.Range("J42").Value2 = "NA"
$ws.Range("D43").Value2 = 557600
$ws.Range("D44").Value2 = 863300
$ws.Range("D45").Value2 = 79000
$ws.Range("D46").Value2 = 2472100
$ws.Range("D47").Value2 = 93700
$ws.Range("E47").Value2 = "NA"
$ws.Range("F47").Value2 = "NA"
$ws.Range("G47").Value2 = "NA"
$ws.Range("H47").Value2 = "NA"
$ws.Range("I47").Value2 = "NA"
$ws.Range("J47").Value2 = "NA"
$ws.Range("D48").Value2 = 585500
$ws.Range("D49").Value2 = "NA"
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 76900
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 3228300
$ws.Range("D57").Value2 = 679600
$ws.Range("D58").Value2 = 8900
$ws.Range("D59").Value2 = 161800
$ws.Range("D60").Value2 = 850200
$ws.Range("D61").Value2 = 88100
$ws.Range("D62").Value2 = 100600
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 1193300
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 1691300
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 2035000
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 301000
$ws.Range("D83").Value2 = 109700
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 568600
$ws.Range("D91").Value2 = -143000
$ws.Range("E91").Value2 = -136000
$ws.Range("G91").Value2 = -118100
$ws.Range("I91").Value2 = -41300
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -319400
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -119700
$ws.Range("D101").Value2 = 6400
$ws.Range("D102").Value2 = 135800
